# Weekly update: insert two new price records (rows 209-210) for
# "Hortaliza, Vega Central Mapocho de Santiago - Pepino ensalada".
# Existing rows 209-222 shift down to 211-224.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 209, pushing the old
# rows 209:222 down to 211:224 (and updating dimension accordingly).
$ws.Range("A209:R210").EntireRow.Insert()

# New row 209 - "Primera" quality record
$ws.Range("A209").Value = 9
$ws.Range("B209").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C209").Value = "Metropolitana"
$ws.Range("D209").Value = 44585
$ws.Range("E209").Value = 13
$ws.Range("F209").Value = 100112043
$ws.Range("G209").Value = "Pepino ensalada"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 70
$ws.Range("K209").Value = 9000
$ws.Range("L209").Value = 10000
$ws.Range("M209").Value = 9500
$ws.Range("N209").Value = "`$/caja 60 unidades"
$ws.Range("O209").Value = "Región de Arica y Parinacota"
$ws.Range("P209").Value = 158
$ws.Range("Q209").Value = 60
$ws.Range("R209").Value = "Hortaliza"

# New row 210 - "Segunda" quality record
$ws.Range("A210").Value = 9
$ws.Range("B210").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C210").Value = "Metropolitana"
$ws.Range("D210").Value = 44585
$ws.Range("E210").Value = 13
$ws.Range("F210").Value = 100112043
$ws.Range("G210").Value = "Pepino ensalada"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Segunda"
$ws.Range("J210").Value = 34
$ws.Range("K210").Value = 7000
$ws.Range("L210").Value = 7000
$ws.Range("M210").Value = 7000
$ws.Range("N210").Value = "`$/caja 100 unidades"
$ws.Range("O210").Value = "Región de Arica y Parinacota"
$ws.Range("P210").Value = 70
$ws.Range("Q210").Value = 100
$ws.Range("R210").Value = "Hortaliza"
